# Replace the "DFPWeiBei-B5" font with "PingFang TC" throughout the
# slide master's title/body default text styles (ppt/slideMasters/slideMaster1.xml).
#
# The OOXML diff changes every <a:latin>, <a:ea>, <a:cs> and <a:sym> typeface
# inside <p:titleStyle> and <p:bodyStyle> (levels 1-9 each) from
# "DFPWeiBei-B5-AZ" to "PingFang TC-AZ".
#
# In the PowerPoint object model these default run-property fonts live on
# Master.TextStyles(ppTitleStyle/ppBodyStyle).Levels(n).Font — the Font
# object's Name / NameFarEast / NameComplexScript / NameOther properties
# correspond to the latin / ea / cs / sym typeface attributes respectively.

$p = $ppt.ActivePresentation
$m = $p.SlideMaster

$oldFont = "DFPWeiBei-B5-AZ"
$newFont = "PingFang TC-AZ"

# ppTitleStyle = 1, ppBodyStyle = 2
$styleIndexes = @(1, 2)

foreach ($styleIdx in $styleIndexes) {
    $ts = $m.TextStyles.Item($styleIdx)

    for ($lvl = 1; $lvl -le 9; $lvl++) {
        $level = $ts.Levels($lvl)
        $font = $level.Font

        $font.Name = $newFont
        $font.NameFarEast = $newFont
        $font.NameComplexScript = $newFont
        $font.NameOther = $newFont
    }
}
